# Commit: "Removed an incorrect underground surface type"
#
# The "Extent of Contamination" worksheet has a table (Table1568) listing,
# per Phase/Category, the fractional breakout of a surface type along with
# its description, distribution parameters, etc. Row 26 was an erroneous
# "Underground" / "HVAC" entry ("The fraction of surface underground which
# is HVAC") that doesn't belong (underground surfaces don't have HVAC) -
# delete that entire table row. Every row below it shifts up by one, the
# table/validation/conditional-formatting ranges shrink by one row, and the
# now-unused shared string gets dropped automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Extent of Contamination")
$ws.Activate()

# Row 26 = Phase "Underground", Category "HVAC" - the erroneous row.
$ws.Rows.Item(26).Delete()

# The conditional formatting rules applied to the table body ("blank
# distribution" / "too many parameters" highlighting) keep their original
# AppliesTo range after the row shift - pull them back in to cover exactly
# the shrunk table body (A2:L33 -> A2:L32).
$fcs = $ws.Cells.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("A2:L32"))
}

# Leave the cursor roughly where it would have ended up after deleting the
# row interactively.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D35").Select()
